$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "22.413.83"
$ws.Range("E2").Value = "  -4.49%  "

# Row 3
$ws.Range("D3").Value = "1.569.99"
$ws.Range("E3").Value = "  -4.68%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("E5").Value = "  +0.00%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.96"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3659"
$ws.Range("E7").Value = "  -3.40%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.38"
$ws.Range("E8").Value = "  -1.56%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3388"
$ws.Range("E9").Value = "  -4.35%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.173"
$ws.Range("E10").Value = "  -3.68%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07586"
$ws.Range("E11").Value = "  -6.19%  "

# Row 12
$ws.Range("E12").Value = "  +0.03%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.19"
$ws.Range("E13").Value = "  -4.13%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.067"
$ws.Range("E14").Value = "  -5.14%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.893"
$ws.Range("E15").Value = "  -6.09%  "

# Row 16
$ws.Range("E16").Value = "  -4.98%  "

# Row 17
$ws.Range("D17").Value = "1.572.12"
$ws.Range("E17").Value = "  -4.81%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.15"
$ws.Range("E18").Value = "  -8.17%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06734"
$ws.Range("E19").Value = "  -3.07%  "

# Row 20
$ws.Range("E20").Value = "  +0.04%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.250"
$ws.Range("E21").Value = "  -7.67%  "

# Row 22
$ws.Range("B22").Value = "BitDAO"
$ws.Range("C22").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.5313"
$ws.Range("E22").Value = "  -7.66%  "

# Row 23
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.44"
$ws.Range("E23").Value = "  -5.49%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.06"
$ws.Range("E24").Value = "  -3.11%  "

# Row 25
$ws.Range("D25").Value = "22.420.92"
$ws.Range("E25").Value = "  -4.50%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.009"
$ws.Range("E26").Value = "  +3.82%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.386"
$ws.Range("E27").Value = "  -4.65%  "

# Row 28
$ws.Range("E28").Value = "  -4.78%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "144.44"
$ws.Range("E29").Value = "  -5.07%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.973"
$ws.Range("E30").Value = "  -4.44%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.43"
$ws.Range("E31").Value = "  -5.65%  "

# Row 32
$ws.Range("D32").Value = "1.749.76"
$ws.Range("E32").Value = "  -4.57%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.054"
$ws.Range("E33").Value = "  +6.43%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.272"
$ws.Range("E34").Value = "  -9.54%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.983"
$ws.Range("E35").Value = "  -7.68%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.39"
$ws.Range("E36").Value = "  -9.60%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02563"
$ws.Range("E37").Value = "  -5.64%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08440"
$ws.Range("E38").Value = "  -3.60%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2304"
$ws.Range("E39").Value = "  -5.54%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06542"
$ws.Range("E40").Value = "  -3.80%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.520"
$ws.Range("E41").Value = "  -7.04%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.87"
$ws.Range("E42").Value = "  -8.80%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.251"
$ws.Range("E43").Value = "  -3.37%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6389"
$ws.Range("E44").Value = "  -7.35%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.58"
$ws.Range("E45").Value = "  -7.46%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9999"
$ws.Range("E46").Value = "  -0.05%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6027"
$ws.Range("E47").Value = "  -5.24%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.782"
$ws.Range("E48").Value = "  -3.27%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.133"
$ws.Range("E49").Value = "  -5.35%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "123.24"
$ws.Range("E50").Value = "  -3.53%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.213"
$ws.Range("E51").Value = "  +2.51%  "
